$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 0.65983699477611935
$ws.Range("C6").Value = 0.13077384523095381
$ws.Range("B7").Value = 0.6578833107588411
$ws.Range("C7").Value = 0.22735452909418116
$ws.Range("B8").Value = 0.65100831550467586
$ws.Range("C8").Value = 0.32393521295740851

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 0.64085424326136353
$ws.Range("C9").Value = 0.42051589682063584

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 0.63298113441377046
$ws.Range("C10").Value = 0.51709658068386322

$ws.Range("A11").Value = 5
$ws.Range("B11").Value = 0.62107385131578929
$ws.Range("C11").Value = 0.6136772645470906

$ws.Range("A12").Value = 6
$ws.Range("B12").Value = 0.6119388
$ws.Range("C12").Value = 0.71025794841031797

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = 0.60260639999999988
$ws.Range("C13").Value = 0.80683863227354535

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = 0.5946072
$ws.Range("C14").Value = 0.90341931613677273

$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 0.58701176275167799
$ws.Range("C15").Value = 1

$ws.Range("D10").Select()
